# StagingTemplates/Staging.IndicatorType.xlsx
#
# The header row on Sheet1 had its "IndicatorType_ID" and "Code" column
# headers swapped (A2 <-> B2); column B/C no longer carry an explicit
# bestFit/custom width override.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Code"
$ws.Range("B2").Value = "IndicatorType_ID"
